$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    2.199488538777095,
    1.926505998055016,
    2.162716768707912,
    2.127502410622359,
    2.251838707917386,
    2.405201540074012,
    1.828924479393044,
    1.685920957457729,
    1.52173442213647,
    1.631068728650424,
    1.655318423829817,
    1.503755554560929,
    2.085209508078397,
    2.065956829980797,
    1.664747774119225,
    1.664587833628723,
    1.443363636114659,
    1.846117816994953,
    1.745270218282298,
    1.695723512637377,
    2.616302171363731
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
